{"js": "// Apply the benchmark-stats update to the single table in the document.\n//\n// Net effect on the (single-column) table:\n//   - rows 0-5 (0-based): simple value replacements\n//   - rows 6-8: deleted (stale raw samples)\n//   - rows 9-11 (now 6-8 after the delete): value replacements\n//   - 3 new rows inserted right after (new) row 8\n//   - the last 3 rows of the table (big multi-tab \"raw dump\" rows) are\n//     each collapsed down to a single summary value\n//\n// Row proxies in this host resolve by *current* collection index at write\n// time, not by stable identity, so after any insert/delete we reload\n// `table.rows` and re-index from the fresh collection before issuing more\n// per-row writes.\n\nconst table = context.document.body.tables.getFirst();\n\n// Helper: replace the full text of a cell while keeping the existing run\n// formatting (rFonts / sz) - a \"Replace\" insert on the cell's range keeps\n// the run's rPr, whereas body.clear() + insertText would drop it.\nfunction setCellText(row, cellIndex, text) {\n  const cell = row.cells.items[cellIndex];\n  const range = cell.body.getRange();\n  range.insertText(text, \"Replace\");\n}\n\n// --- Step 1: rows 0-5, plain value swaps -----------------------------\ntable.rows.load(\"items\");\nawait context.sync();\nlet rows = table.rows.items;\n\nsetCellText(rows[0], 0, \"0M\");\nsetCellText(rows[1], 0, \"0M\");\nsetCellText(rows[2], 0, \"0M\");\nsetCellText(rows[3], 0, \"274\");\nsetCellText(rows[4], 0, \"0.00002\");\nsetCellText(rows[5], 0, \"0.00066\");\nawait context.sync();\n\n// --- Step 2: delete the three stale sample rows (old rows 6,7,8) -----\ntable.rows.load(\"items\");\nawait context.sync();\nrows = table.rows.items;\n\n// delete from the bottom up so earlier indices stay valid within this batch\nrows[8].delete();\nrows[7].delete();\nrows[6].delete();\nawait context.sync();\n\n// --- Step 3: update the rows that shifted into slots 6-8 -------------\ntable.rows.load(\"items\");\nawait context.sync();\nrows = table.rows.items;\n\nsetCellText(rows[6], 0, \"0.00018\"); // unchanged value, rewritten for clarity\nsetCellText(rows[7], 0, \"0.00005\");\nsetCellText(rows[8], 0, \"0.00028\");\nawait context.sync();\n\n// --- Step 4: insert three new rows after (new) row 8 ------------------\ntable.rows.load(\"items\");\nawait context.sync();\nrows = table.rows.items;\n\nrows[8].insertRows(\"After\", 3, [[\"0.00034\"], [\"0.00042\"], [\"0.05624\"]]);\nawait context.sync();\n\n// --- Step 5: collapse the last three (multi-tab) rows to single values\n// Net row count is unchanged (3 deleted + 3 inserted), so these are still\n// the final three rows of the table.\ntable.rows.load(\"items\");\nawait context.sync();\nrows = table.rows.items;\nconst n = rows.length;\n\nsetCellText(rows[n - 3], 0, \"99.92\");\nsetCellText(rows[n - 2], 0, \"0.06\");\nsetCellText(rows[n - 1], 0, \"67\");\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-stats update to the single table in the document.\n#\n# Net effect on the (single-column) table:\n#   - rows 1-6 (1-based): simple value replacements\n#   - rows 7-9: deleted (stale raw samples: 0.00014 / 0.00004 / 0.00016)\n#   - rows 7-9 (after the delete, now holding the old rows 10-12):\n#     value replacements\n#   - 3 new rows inserted right after (new) row 9\n#   - the last 3 rows of the table (big multi-tab \"raw dump\" rows) are\n#     each collapsed down to a single summary value\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Step 1: rows 1-6, plain value swaps ------------------------------\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"274\"\n$t.Cell(5, 1).Range.Text = \"0.00002\"\n$t.Cell(6, 1).Range.Text = \"0.00066\"\n\n# --- Step 2: delete the three stale sample rows (old rows 7,8,9) -----\n# delete from the bottom up so earlier indices stay valid within this batch\n$t.Rows.Item(9).Delete()\n$t.Rows.Item(8).Delete()\n$t.Rows.Item(7).Delete()\n\n# --- Step 3: update the rows that shifted into slots 7-9 --------------\n$t.Cell(7, 1).Range.Text = \"0.00018\"  # unchanged value, rewritten for clarity\n$t.Cell(8, 1).Range.Text = \"0.00005\"\n$t.Cell(9, 1).Range.Text = \"0.00028\"\n\n# --- Step 4: insert three new rows after (new) row 9 -------------------\n$refRow = $t.Rows.Item(10)\n$t.Rows.Add($refRow) | Out-Null\n$t.Rows.Add($refRow) | Out-Null\n$t.Rows.Add($refRow) | Out-Null\n$t.Cell(10, 1).Range.Text = \"0.00034\"\n$t.Cell(11, 1).Range.Text = \"0.00042\"\n$t.Cell(12, 1).Range.Text = \"0.05624\"\n\n# --- Step 5: collapse the last three (multi-tab) rows to single values\n# Net row count is unchanged (3 deleted + 3 inserted), so these are still\n# the final three rows of the table.\n$n = $t.Rows.Count\n$t.Cell($n - 2, 1).Range.Text = \"99.92\"\n$t.Cell($n - 1, 1).Range.Text = \"0.06\"\n$t.Cell($n, 1).Range.Text = \"67\"\n"}
